# Append: 2025-11-23 01:38 JST
#
# A fresh scrape ran and produced one brand-new listing
# (https://www.lancers.jp/work/detail/5439670), which the scraper inserts
# right after the header row. Every previously-seen row shifts down by one
# and gets its "fetched at" timestamp refreshed to the new run time; the
# previously-newest row (5439402) simply slides from row 5 to row 6.
#
# We rebuild the sheet explicitly (rather than relying on Rows.Insert, whose
# hyperlink-shifting semantics don't line up with the desired end state) so
# every cell value/hyperlink lands exactly where the new snapshot expects it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$timestamp = "2025-11-23 01:38:33"

# Drop the old hyperlinks up front; we'll re-add one per F-column URL cell
# below once all the row data is in its final place, so ref/target stay in
# sync (no stale anchors left pointing at shifted cells).
if ($ws.Hyperlinks.Count -gt 0) {
    $ws.Hyperlinks.Delete()
}

# --- Row 2 (unchanged listing, refreshed timestamp only) ---
$ws.Range("A2").Value = $timestamp
$ws.Range("B2").Value = "【技術者募集】家庭用消臭デバイス「Maneki Air」開発"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5439445"
$ws.Range("G2").Value = 368
$ws.Range("H2").Value = "🔥AI,Ai ◆開発"

# --- Row 3 (unchanged listing, refreshed timestamp only) ---
$ws.Range("A3").Value = $timestamp
$ws.Range("B3").Value = "【急募】掲示板サイト(爆サイ)自動書き込みソフト開発者募集"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5439484"
$ws.Range("G3").Value = 93
$ws.Range("H3").Value = "◆開発 ◇サイト"

# --- Row 4 (unchanged listing, refreshed timestamp only) ---
$ws.Range("A4").Value = $timestamp
$ws.Range("B4").Value = "名刺/プロフィール共有アプリ開発"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "1,000 ~ 5,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5439373"
$ws.Range("G4").Value = 85
$ws.Range("H4").Value = "◆開発 ◇アプリ"

# --- Row 5 (NEW listing from this run) ---
$ws.Range("A5").Value = $timestamp
$ws.Range("B5").Value = "【緊急】海外からWordPress管理画面にログインできない問題の調査と修正"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5439670"
$ws.Range("G5").Value = 50
$ws.Range("H5").Value = "◇管理 ○WordPress"

# --- Row 6 (was row 5 before this run; same listing, shifted down) ---
$ws.Range("A6").Value = $timestamp
$ws.Range("B6").Value = "【緊急】海外からWordPress管理画面にログインできない問題の調査と修正"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5439402"
$ws.Range("G6").Value = 50
$ws.Range("H6").Value = "◇管理 ○WordPress"

# --- Row 7 (was row 6 before this run; same listing, shifted down) ---
$ws.Range("A7").Value = $timestamp
$ws.Range("B7").Value = "限定公開 限定公開の仕事"
$ws.Range("C7").Value = "システム開発"
$ws.Range("D7").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E7").Value = "期限情報なし"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5439488"
$ws.Range("G7").Value = 25

# --- Row 8 (was row 7 before this run; same listing, shifted down) ---
$ws.Range("A8").Value = $timestamp
$ws.Range("B8").Value = "【急募】ファン応援プラットフォームの構築をお手伝いください!"
$ws.Range("C8").Value = "システム開発"
$ws.Range("D8").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E8").Value = "期限情報なし"
$ws.Range("F8").Value = "https://www.lancers.jp/work/detail/5439395"
$ws.Range("G8").Value = 18

# --- Row 9 (was row 8 before this run; same listing, shifted down) ---
$ws.Range("A9").Value = $timestamp
$ws.Range("B9").Value = "【Stable Diffusion】参考動画に沿って約100プロンプト構築"
$ws.Range("C9").Value = "システム開発"
$ws.Range("D9").Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Range("E9").Value = "期限情報なし"
$ws.Range("F9").Value = "https://www.lancers.jp/work/detail/5432055"
$ws.Range("G9").Value = 10

# Re-create the F-column hyperlinks (URL text doubles as the display text,
# matching the scraper's original output) now that every row sits in its
# final position.
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5439445")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5439484")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5439373")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5439670")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5439402")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5439488")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5439395")
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.lancers.jp/work/detail/5432055")
